$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 526.4091
$ws.Range("J17").Value = 526.4091
$ws.Range("L17").Value = 1579.2273
$ws.Range("N17").Value = -1915.2273

# Row 51
$ws.Range("H51").Value = 5737.421
$ws.Range("J51").Value = 6437.9287
$ws.Range("L51").Value = 6437.9287
$ws.Range("N51").Value = -7405.9287

# Row 105
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

# Row 113
$ws.Range("H113").Value = 3307.05
$ws.Range("I113").Value = 3042.7334
$ws.Range("K113").Value = 3042.7334
$ws.Range("M113").Value = 211.2665999999999

# Row 127
$ws.Range("H127").Value = 847
$ws.Range("I127").Value = 569.4
$ws.Range("J127").Value = 1772.3334
$ws.Range("K127").Value = 1708.2
$ws.Range("L127").Value = 5317.0002
$ws.Range("M127").Value = 3251.8
$ws.Range("N127").Value = -15237.0002

# Row 137
$ws.Range("H137").Value = 2705428.2
$ws.Range("I137").Value = 3705840.5
$ws.Range("J137").Value = 4314.8
$ws.Range("K137").Value = 11117521.5
$ws.Range("L137").Value = 12944.4
$ws.Range("M137").Value = -11114971.5
$ws.Range("N137").Value = -18044.4

# Row 138
$ws.Range("H138").Value = 3273487
$ws.Range("I138").Value = 234186.8
$ws.Range("J138").Value = 15154388
$ws.Range("K138").Value = 702560.3999999999
$ws.Range("L138").Value = 45463164
$ws.Range("M138").Value = -697420.3999999999
$ws.Range("N138").Value = -45473444

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 23451.666
$ws.Range("J24").Value = 23451.666
$ws.Range("L24").Value = 23451.666
$ws.Range("N24").Value = -24199.666

# Row 32
$ws.Range("H32").Value = 1265.58
$ws.Range("I32").Value = 1038.9186
$ws.Range("J32").Value = 2657.9285
$ws.Range("K32").Value = 1038.9186
$ws.Range("L32").Value = 2657.9285
$ws.Range("M32").Value = -751.9186
$ws.Range("N32").Value = -3231.9285

# Row 74
$ws.Range("H74").Value = 5041893.5
$ws.Range("I74").Value = 5977499.5
$ws.Range("J74").Value = 129964.25
$ws.Range("K74").Value = 5977499.5
$ws.Range("L74").Value = 129964.25
$ws.Range("M74").Value = -5976625.5
$ws.Range("N74").Value = -131712.25

# Row 77
$ws.Range("H77").Value = 5041893.5
$ws.Range("I77").Value = 5977499.5
$ws.Range("J77").Value = 129964.25
$ws.Range("K77").Value = 29887497.5
$ws.Range("L77").Value = 649821.25
$ws.Range("M77").Value = -29883129.5
$ws.Range("N77").Value = -658557.25

# Row 100
$ws.Range("H100").Value = 23451.666
$ws.Range("J100").Value = 23451.666
$ws.Range("L100").Value = 23451.666
$ws.Range("N100").Value = -25615.666

# Row 132
$ws.Range("H132").Value = 73159.73
$ws.Range("I132").Value = 44400.824
$ws.Range("J132").Value = 183402.17
$ws.Range("K132").Value = 133202.472
$ws.Range("L132").Value = 550206.51
$ws.Range("M132").Value = -130672.472
$ws.Range("N132").Value = -555266.51

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 250000800
$ws.Range("I22").Value = 333334000
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 333334000
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -333333650
$ws.Range("N22").Value = -1900

# Row 62
$ws.Range("H62").Value = 2925
$ws.Range("I62").Value = 2850
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2850
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2226
$ws.Range("N62").Value = -4248

# Row 65
$ws.Range("H65").Value = 2925
$ws.Range("I65").Value = 2850
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14250
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11130
$ws.Range("N65").Value = -21240

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 620.9091
$ws.Range("J15").Value = 812.5
$ws.Range("L15").Value = 2437.5
$ws.Range("N15").Value = -2717.5

# Row 122
$ws.Range("H122").Value = 723.4091
$ws.Range("I122").Value = 274.66666
$ws.Range("J122").Value = 1261.9
$ws.Range("K122").Value = 2471.99994
$ws.Range("L122").Value = 11357.1
$ws.Range("M122").Value = -21.9999399999997
$ws.Range("N122").Value = -16257.1

# Row 131
$ws.Range("H131").Value = 14926309
$ws.Range("I131").Value = 83333710
$ws.Range("J131").Value = 1057.5818
$ws.Range("K131").Value = 250001130
$ws.Range("L131").Value = 3172.7454
$ws.Range("M131").Value = -249996090
$ws.Range("N131").Value = -13252.7454

# Row 137
$ws.Range("H137").Value = 26458.36
$ws.Range("I137").Value = 1950
$ws.Range("J137").Value = 31126.62
$ws.Range("K137").Value = 5850
$ws.Range("L137").Value = 93379.86
$ws.Range("M137").Value = -750
$ws.Range("N137").Value = -103579.86

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 36812.87
$ws.Range("I70").Value = 46033.332
$ws.Range("J70").Value = 5199.857
$ws.Range("K70").Value = 46033.332
$ws.Range("L70").Value = 5199.857
$ws.Range("M70").Value = -45763.332
$ws.Range("N70").Value = -5739.857

# Row 73
$ws.Range("H73").Value = 36812.87
$ws.Range("I73").Value = 46033.332
$ws.Range("J73").Value = 5199.857
$ws.Range("K73").Value = 46033.332
$ws.Range("L73").Value = 5199.857
$ws.Range("M73").Value = -45097.332
$ws.Range("N73").Value = -7071.857

# Row 101
$ws.Range("H101").Value = 44996.668
$ws.Range("J101").Value = 44996.668
$ws.Range("L101").Value = 44996.668
$ws.Range("N101").Value = -51486.668

# Row 122
$ws.Range("H122").Value = 2907.0625
$ws.Range("I122").Value = 2581.1538
$ws.Range("J122").Value = 4319.3335
$ws.Range("K122").Value = 7743.4614
$ws.Range("L122").Value = 12958.0005
$ws.Range("M122").Value = -5293.4614
$ws.Range("N122").Value = -17858.0005

# Row 141
$ws.Range("H141").Value = 31272.5
$ws.Range("J141").Value = 31272.5
$ws.Range("L141").Value = 31272.5
$ws.Range("N141").Value = -41632.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 918.7222
$ws.Range("I22").Value = 591.375
$ws.Range("J22").Value = 1180.6
$ws.Range("K22").Value = 591.375
$ws.Range("L22").Value = 1180.6
$ws.Range("M22").Value = -296.375
$ws.Range("N22").Value = -1770.6

# Row 27
$ws.Range("H27").Value = 918.7222
$ws.Range("I27").Value = 591.375
$ws.Range("J27").Value = 1180.6
$ws.Range("K27").Value = 591.375
$ws.Range("L27").Value = 1180.6
$ws.Range("M27").Value = -484.375
$ws.Range("N27").Value = -1394.6

# Row 68
$ws.Range("H68").Value = 1669.3334
$ws.Range("I68").Value = 1622.0952
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1622.0952
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -873.0952
$ws.Range("N68").Value = -3498

# Row 71
$ws.Range("H71").Value = 1669.3334
$ws.Range("I71").Value = 1622.0952
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 8110.476
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -4366.476
$ws.Range("N71").Value = -17488

$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("H104").Value = 27500
$ws.Range("J104").Value = 27500
$ws.Range("L104").Value = 27500
$ws.Range("N104").Value = -34488

# Row 113
$ws.Range("H113").Value = 1063.12
$ws.Range("I113").Value = 755.7857
$ws.Range("J113").Value = 1454.2727
$ws.Range("K113").Value = 2267.3571
$ws.Range("L113").Value = 4362.8181
$ws.Range("M113").Value = -97.35710000000017
$ws.Range("N113").Value = -8702.8181

# Row 129
$ws.Range("H129").Value = 32826.668
$ws.Range("J129").Value = 32826.668
$ws.Range("L129").Value = 32826.668
$ws.Range("N129").Value = -42826.668

# Row 136
$ws.Range("H136").Value = 51483.4
$ws.Range("I136").Value = 42675.25
$ws.Range("J136").Value = 64695.625
$ws.Range("K136").Value = 128025.75
$ws.Range("L136").Value = 194086.875
$ws.Range("M136").Value = -125475.75
$ws.Range("N136").Value = -199186.875

# Row 140
$ws.Range("H140").Value = 40582.715
$ws.Range("J140").Value = 40582.715
$ws.Range("L140").Value = 40582.715
$ws.Range("N140").Value = -50942.715

Write-Output "Applied all Hades_Profits updates"
